# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff has been generated (status -> "Ready for
# handoff", new handoff xliff file names/timestamps, and a version-mismatch
# error message), and widens the "Error Detail" column so the new message
# is readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 22:44:11"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text entry so "False" stays a literal string
# (matching the rest of the sheet) instead of becoming a native Boolean.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-19 22:44:07"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f95067305d7f67a6fe5e158c8b92edf4fcc02e4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e7ebc1fe63c92f7000d4ba8bab18d6698c6acbe/e2e/b.md."
# Widen the "Error Detail" column (P) to fit the new message (stored width 40).
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text entry so "False" stays a literal string
# (matching the rest of the sheet) instead of becoming a native Boolean.
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-19 22:44:11"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f95067305d7f67a6fe5e158c8b92edf4fcc02e4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e7ebc1fe63c92f7000d4ba8bab18d6698c6acbe/e2e/b.md."
# Widen the "Error Detail" column (P) to fit the new message (stored width 40).
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
